$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the existing EMPLOYEES row (row 3 had a stray "EMPLOYEES_NEW" name) ---
$ws.Range("A3").Value = "EMPLOYEES"

# --- Add the two new data rows -------------------------------------------
# Values are written in the same order the strings first appear in the
# author's shared-string table, so that the underlying shared-string pool
# is rebuilt in the same sequence.
$ws.Cells.Item(5, 3).Value = "select * from stg_department"
$ws.Cells.Item(4, 3).Value = "select * from stg_employees"
$ws.Cells.Item(5, 4).Value = "Department_ID"
$ws.Cells.Item(4, 1).Value = "STG_EMPLOYEES"
$ws.Cells.Item(5, 1).Value = "DEPARTMENTS"
$ws.Cells.Item(5, 2).Value = "SELECT * FROM HR.DEPARTMENTS"

$ws.Cells.Item(4, 2).Value = "SELECT * FROM HR.EMPLOYEES"
$ws.Cells.Item(4, 4).Value = "EMPLOYEE_ID"
$ws.Cells.Item(4, 5).Value = "Y"
$ws.Cells.Item(5, 5).Value = "Y"

# New rows should carry the same 10pt font as the rest of the table.
$ws.Range("A4:E5").Font.Size = 10

# --- Formatting clean-up ---------------------------------------------------
# The table no longer needs the tall, word-wrapped rows - turn wrap text
# off for every data row and let the row heights shrink back to normal.
$ws.Range("A2:E5").WrapText = $false
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# Columns A:C must grow to fit the new, longer text values.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# Leave the selection where the author left it when they saved.
$ws.Range("C9").Select() | Out-Null
